$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for year 2021 (row 12), following the same layout as row 11 (2020)
$targetRow = 12
$sourceRow = 11

# Copy formatting (style) from the previous data row so the new row matches
# the existing look (bold/centered/bordered year label in column A).
$ws.Range("A" + $sourceRow + ":V" + $sourceRow).Copy() | Out-Null
$ws.Range("A" + $targetRow + ":V" + $targetRow).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Year label
$ws.Cells.Item($targetRow, 1).Value = "2021年"

# Numeric columns with values
$ws.Cells.Item($targetRow, 2).Value = 152899.9    # B
$ws.Cells.Item($targetRow, 5).Value = 1827786.4   # E
$ws.Cells.Item($targetRow, 6).Value = 1745261.3   # F
$ws.Cells.Item($targetRow, 7).Value = 2047577.7   # G
$ws.Cells.Item($targetRow, 8).Value = 1879744.1   # H
$ws.Cells.Item($targetRow, 9).Value = 755305.7    # I
$ws.Cells.Item($targetRow, 11).Value = 2191700.1  # K
$ws.Cells.Item($targetRow, 12).Value = 1591918.1  # L
$ws.Cells.Item($targetRow, 13).Value = 2038356.9  # M
$ws.Cells.Item($targetRow, 18).Value = 827834.9   # R
$ws.Cells.Item($targetRow, 20).Value = 3783618.2  # T
$ws.Cells.Item($targetRow, 22).Value = 76087.7    # V

# Columns with no data for this year remain as empty strings (matching the
# existing rows' representation of "no value" cells)
$ws.Cells.Item($targetRow, 3).Value = ""   # C
$ws.Cells.Item($targetRow, 4).Value = ""   # D
$ws.Cells.Item($targetRow, 10).Value = ""  # J
$ws.Cells.Item($targetRow, 14).Value = ""  # N
$ws.Cells.Item($targetRow, 15).Value = ""  # O
$ws.Cells.Item($targetRow, 16).Value = ""  # P
$ws.Cells.Item($targetRow, 17).Value = ""  # Q
$ws.Cells.Item($targetRow, 19).Value = ""  # S
$ws.Cells.Item($targetRow, 21).Value = ""  # U
